# EI Variable Installments T2 scenarios
#
# Adds a new "waittopageload1" automation step as row 6 of the
# "Edit Repayment Schedule" sheet (pushing the previous rows 6-12 down
# to 7-13), and makes "Edit Repayment Schedule" the active sheet/tab
# (previously "Summary" was active).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Edit Repayment Schedule")

# Insert a blank row above the old row 6 ("clickonsubmit" / "Submit"),
# shifting it - and everything below it - down by one.
$ws.Rows.Item(6).Insert()

# Give the new row 6 the same look as the other "wait" row (row 3:
# "waittopageload" / 2000) and fill in its content.
$ws.Range("A3:B3").Copy($ws.Range("A6:B6"))
$excel.CutCopyMode = $false

$ws.Range("A6").Value = "waittopageload1"
$ws.Range("B6").Value = 2000

# Make this sheet the active tab, with the new row selected.
$ws.Activate()
$ws.Range("A6:B6").Select()
